# Updates the cryptocurrency price/volume table (columns D "Price" and
# E "Volume(1h)") on Sheet1 with refreshed values, per the
# "Updated cryptos list ... with GitHub Actions" commit.
#
# Values are written with a leading apostrophe so Excel stores them as
# text (matching the workbook's original inlineStr cells) instead of
# silently reinterpreting number-looking strings (e.g. "0.9988") as
# numeric values or percent-looking strings (e.g. "+0.39%") as
# percentages. ClearFormats() afterwards strips the quote-prefix
# formatting Excel applies for the apostrophe trick, restoring the
# cells' original (default/general) style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''29.233.62'
$ws.Range("E2").Value = '''  +0.39%  '
$ws.Range("D3").Value = '''1.844.62'
$ws.Range("E3").Value = '''  +0.67%  '
$ws.Range("D4").Value = '''0.9988'
$ws.Range("E4").Value = '''  -0.04%  '
$ws.Range("D5").Value = '''240.65'
$ws.Range("E5").Value = '''  +0.07%  '
$ws.Range("D6").Value = '''0.6727'
$ws.Range("E6").Value = '''  -1.60%  '
$ws.Range("D7").Value = '''0.9997'
$ws.Range("E7").Value = '''  -0.02%  '
$ws.Range("D8").Value = '''0.07431'
$ws.Range("E8").Value = '''  -0.18%  '
$ws.Range("D9").Value = '''0.2954'
$ws.Range("E9").Value = '''  -1.89%  '
$ws.Range("D10").Value = '''22.89'
$ws.Range("E10").Value = '''  -0.77%  '
$ws.Range("D11").Value = '''0.07717'
$ws.Range("E11").Value = '''  +0.76%  '
$ws.Range("D12").Value = '''1.837.17'
$ws.Range("E12").Value = '''  +0.21%  '
$ws.Range("D13").Value = '''5.008'
$ws.Range("E13").Value = '''  -0.84%  '
$ws.Range("D14").Value = '''0.6730'
$ws.Range("E14").Value = '''  -1.15%  '
$ws.Range("D15").Value = '''86.18'
$ws.Range("E15").Value = '''  -1.48%  '
$ws.Range("D16").Value = '''6.123'
$ws.Range("E16").Value = '''  -0.26%  '
$ws.Range("D17").Value = '''0.000008316'
$ws.Range("E17").Value = '''  +1.83%  '
$ws.Range("D18").Value = '''29.155.24'
$ws.Range("E18").Value = '''  +0.16%  '
$ws.Range("D19").Value = '''228.74'
$ws.Range("E20").Value = '''  +0.18%  '
$ws.Range("D21").Value = '''1.001'
$ws.Range("E21").Value = '''  +0.07%  '
$ws.Range("D22").Value = '''7.189'
$ws.Range("E22").Value = '''  -2.93%  '
$ws.Range("D23").Value = '''1.0000'
$ws.Range("E23").Value = '''  +0.03%  '
$ws.Range("D24").Value = '''160.60'
$ws.Range("E24").Value = '''  +0.38%  '
$ws.Range("D25").Value = '''8.692'
$ws.Range("E25").Value = '''  -0.48%  '
$ws.Range("D26").Value = '''0.1404'
$ws.Range("E26").Value = '''  -3.44%  '
$ws.Range("E27").Value = '''  -0.48%  '
$ws.Range("D28").Value = '''1.507'
$ws.Range("E28").Value = '''  -0.14%  '
$ws.Range("D29").Value = '''4.180'
$ws.Range("E29").Value = '''  -2.14%  '
$ws.Range("D30").Value = '''4.073'
$ws.Range("E30").Value = '''  -1.82%  '
$ws.Range("D31").Value = '''1.190'
$ws.Range("E31").Value = '''  -0.85%  '
$ws.Range("D32").Value = '''0.05306'
$ws.Range("E32").Value = '''  +2.83%  '
$ws.Range("E33").Value = '''  +2.16%  '
$ws.Range("D34").Value = '''0.7582'
$ws.Range("E34").Value = '''  -0.97%  '
$ws.Range("E35").Value = '''  +0.35%  '
$ws.Range("D36").Value = '''2.675'
$ws.Range("E36").Value = '''  +0.06%  '
$ws.Range("D37").Value = '''1.326.79'
$ws.Range("E37").Value = '''  +1.42%  '
$ws.Range("D38").Value = '''0.01805'
$ws.Range("E38").Value = '''  -1.52%  '
$ws.Range("E39").Value = '''  +0.63%  '
$ws.Range("D40").Value = '''0.9190'
$ws.Range("E40").Value = '''  -1.24%  '
$ws.Range("D41").Value = '''5.954'
$ws.Range("E41").Value = '''  +2.58%  '
$ws.Range("D42").Value = '''1.001'
$ws.Range("D43").Value = '''103.42'
$ws.Range("E43").Value = '''  -0.93%  '
$ws.Range("D44").Value = '''0.08034'
$ws.Range("E44").Value = '''  +15.86%  '
$ws.Range("D45").Value = '''1.970.61'
$ws.Range("E45").Value = '''  -0.55%  '
$ws.Range("D46").Value = '''0.5161'
$ws.Range("E46").Value = '''  -0.62%  '
$ws.Range("E47").Value = '''  +0.36%  '
$ws.Range("E48").Value = '''  -2.12%  '
$ws.Range("D49").Value = '''63.90'
$ws.Range("E49").Value = '''  -2.03%  '
$ws.Range("D50").Value = '''9.201'
$ws.Range("E50").Value = '''  -3.49%  '
$ws.Range("D51").Value = '''0.05949'

# Restore default (General) formatting on the edited cells; the
# apostrophe-as-text trick above stamps a transient quote-prefix style
# that isn't present in the source workbook.
$ws.Range("D2:E51").ClearFormats()
